$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 24/25: Toncoin and WrappedliquidstakedEther2.0 swap positions with new data ---
$ws.Range("B24").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C24").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D24").Value = '2.186.41'
$ws.Range("E24").Value = '  +2.79%  '

$ws.Range("B25").Value = 'Toncoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.119'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.96%  '

# --- Remaining price/volume updates ---
$ws.Range("D2").Value = '28.643.05'
$ws.Range("E2").Value = '  +2.14%  '
$ws.Range("D3").Value = '1.871.15'
$ws.Range("E3").Value = '  +2.22%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.006'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.53%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '325.88'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.13%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.004'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.37%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4645'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.27%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3884'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.26%  '
$ws.Range("E9").Value = '  +0.18%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9750'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.52%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '21.98'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.46%  '
$ws.Range("D12").Value = '1.894.65'
$ws.Range("E12").Value = '  -1.44%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.995'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.46%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.703'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.66%  '
$ws.Range("E15").Value = '  +3.47%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '88.02'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.06%  '
$ws.Range("E17").Value = '  +0.39%  '
$ws.Range("E18").Value = '  +1.08%  '
$ws.Range("E19").Value = '  +1.28%  '
$ws.Range("E20").Value = '  +0.29%  '
$ws.Range("D21").Value = '28.660.33'
$ws.Range("E21").Value = '  +2.12%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.287'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.40%  '
$ws.Range("E23").Value = '  +0.11%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '152.81'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.57%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.21'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.36%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.790'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.43%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.987'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.62%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '119.36'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.90%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09363'
$ws.Range("D31").Style = "Normal"
$ws.Range("E32").Value = '  -2.11%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.266'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.54%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.337'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.62%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.330'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.49%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.05796'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.35%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02102'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.95%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.147'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.32%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '7.766'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.18%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5628'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.62%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1786'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.39%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '9.769'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.29%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.07207'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.64%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '11.76'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.96%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.5313'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.71%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.149'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.06%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.106'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.83%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.826'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.30%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '113.20'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.06%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.423'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +4.56%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.003'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.39%  '
